$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.554.88"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "3.093.17"
$ws.Range("E3").Value = "  -2.66%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.06"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.45"
$ws.Range("E6").Value = "  -4.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.091.47"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.467"
$ws.Range("E9").Value = "  +3.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  -3.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.404"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "3.624.45"
$ws.Range("E14").Value = "  -2.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.38"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("E16").Value = "  -3.59%  "
$ws.Range("D17").Value = "57.583.50"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "3.083.94"
$ws.Range("E18").Value = "  -2.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.91"
$ws.Range("E19").Value = "  -4.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.53"
$ws.Range("E20").Value = "  -3.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.90"
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "349.76"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.60"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "0.0₃0864"
$ws.Range("E28").Value = "  -10.01%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.13"
$ws.Range("E30").Value = "  -5.92%  "
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.99"
$ws.Range("E32").Value = "  -9.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.14"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.86"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.05"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("E36").Value = "  -7.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.02"
$ws.Range("E37").Value = "  -4.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.63"
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.24"
$ws.Range("E39").Value = "  -5.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0663"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.692"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").Value = "2.391.33"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.90"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "3.132.55"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("E48").Value = "  -4.24%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.02"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.952"
$ws.Range("E50").Value = "  -7.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.50"
$ws.Range("E51").Value = "  -6.11%  "
